$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 86
$ws.Range("H86").Value = 37056884
$ws.Range("I86").Value = 3162.1875
$ws.Range("J86").Value = 90953200
$ws.Range("K86").Value = 3162.1875
$ws.Range("L86").Value = 90953200
$ws.Range("M86").Value = -2039.1875
$ws.Range("N86").Value = -90955446
# row 89
$ws.Range("H89").Value = 37056884
$ws.Range("I89").Value = 3162.1875
$ws.Range("J89").Value = 90953200
$ws.Range("K89").Value = 15810.9375
$ws.Range("L89").Value = 454766000
$ws.Range("M89").Value = -10194.9375
$ws.Range("N89").Value = -454777232
# row 138
$ws.Range("H138").Value = 3663.7144
$ws.Range("J138").Value = 4771.1177
$ws.Range("L138").Value = 14313.3531
$ws.Range("N138").Value = -24593.3531

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 35191.234
$ws.Range("I2").Value = 1380.8422
$ws.Range("K2").Value = 1380.8422
$ws.Range("M2").Value = -1267.8422
# row 32
$ws.Range("H32").Value = 4260.1816
$ws.Range("I32").Value = 4260.1816
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4260.1816
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3973.1816
$ws.Range("N32").ClearContents()
# row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# row 116
$ws.Range("H116").Value = 35191.234
$ws.Range("I116").Value = 1380.8422
$ws.Range("K116").Value = 1380.8422
$ws.Range("M116").Value = 913.1578
# row 122
$ws.Range("H122").Value = 1472855.5
$ws.Range("I122").Value = 6387.4614
$ws.Range("J122").Value = 2594272.2
$ws.Range("K122").Value = 19162.3842
$ws.Range("L122").Value = 7782816.600000001
$ws.Range("M122").Value = -16712.3842
$ws.Range("N122").Value = -7787716.600000001
# row 138
$ws.Range("H138").Value = 122887.664
$ws.Range("J138").Value = 122887.664
$ws.Range("L138").Value = 122887.664
$ws.Range("N138").Value = -133167.664

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 35191.234
$ws.Range("I3").Value = 1380.8422
$ws.Range("K3").Value = 1380.8422
$ws.Range("M3").Value = -1266.8422
# row 94
$ws.Range("H94").Value = 50667.5
$ws.Range("I94").Value = 92335.664
$ws.Range("J94").Value = 8999.333000000001
$ws.Range("K94").Value = 92335.664
$ws.Range("L94").Value = 8999.333000000001
$ws.Range("M94").Value = -91884.664
$ws.Range("N94").Value = -9901.333000000001
# row 107
$ws.Range("H107").Value = 2463.1667
$ws.Range("I107").Value = 2444.3157
$ws.Range("J107").Value = 2534.8
$ws.Range("K107").Value = 2444.3157
$ws.Range("L107").Value = 2534.8
$ws.Range("M107").Value = -524.3157000000001
$ws.Range("N107").Value = -6374.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 2154
$ws.Range("I16").Value = 870.1539
$ws.Range("J16").Value = 10499
$ws.Range("K16").Value = 870.1539
$ws.Range("L16").Value = 10499
$ws.Range("M16").Value = -583.1539
$ws.Range("N16").Value = -11073
# row 31
$ws.Range("H31").Value = 4389.952
$ws.Range("I31").Value = 3452.842
$ws.Range("J31").Value = 5164.087
$ws.Range("K31").Value = 3452.842
$ws.Range("L31").Value = 5164.087
$ws.Range("M31").Value = -3157.842
$ws.Range("N31").Value = -5754.087
# row 34
$ws.Range("H34").Value = 4389.952
$ws.Range("I34").Value = 3452.842
$ws.Range("J34").Value = 5164.087
$ws.Range("K34").Value = 3452.842
$ws.Range("L34").Value = 5164.087
$ws.Range("M34").Value = -3250.842
$ws.Range("N34").Value = -5568.087
# row 107
$ws.Range("H107").Value = 41677836
$ws.Range("J107").Value = 3383.9167
$ws.Range("L107").Value = 3383.9167
$ws.Range("N107").Value = -7223.9167
# row 113
$ws.Range("H113").Value = 2154
$ws.Range("I113").Value = 870.1539
$ws.Range("J113").Value = 10499
$ws.Range("K113").Value = 870.1539
$ws.Range("L113").Value = 10499
$ws.Range("M113").Value = 1299.8461
$ws.Range("N113").Value = -14839
# row 114
$ws.Range("H114").Value = 75021.336
$ws.Range("I114").Value = 50621
$ws.Range("J114").Value = 87221.5
$ws.Range("K114").Value = 50621
$ws.Range("L114").Value = 87221.5
$ws.Range("M114").Value = -46282
$ws.Range("N114").Value = -95899.5
# row 118
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 68
$ws.Range("H68").Value = 25007272
$ws.Range("I68").Value = 1067.625
$ws.Range("J68").Value = 41678076
$ws.Range("K68").Value = 3202.875
$ws.Range("L68").Value = 125034228
$ws.Range("M68").Value = -2391.875
$ws.Range("N68").Value = -125035850
# row 71
$ws.Range("H71").Value = 25007272
$ws.Range("I71").Value = 1067.625
$ws.Range("J71").Value = 41678076
$ws.Range("K71").Value = 9608.625
$ws.Range("L71").Value = 375102684
$ws.Range("M71").Value = -5552.625
$ws.Range("N71").Value = -375110796
# row 107
$ws.Range("H107").Value = 673.13635
$ws.Range("I107").Value = 276.9091
$ws.Range("J107").Value = 805.2121
$ws.Range("K107").Value = 830.7273
$ws.Range("L107").Value = 2415.6363
$ws.Range("M107").Value = 1089.2727
$ws.Range("N107").Value = -6255.6363
# row 117
$ws.Range("H117").Value = 958.9231
$ws.Range("I117").Value = 980.5
$ws.Range("J117").Value = 700
$ws.Range("K117").Value = 2941.5
$ws.Range("L117").Value = 2100
$ws.Range("M117").Value = 500.5
$ws.Range("N117").Value = -8984
# row 121
$ws.Range("H121").Value = 1177654.9
$ws.Range("I121").Value = 716.875
$ws.Range("J121").Value = 2223822
$ws.Range("K121").Value = 2150.625
$ws.Range("L121").Value = 6671466
$ws.Range("M121").Value = -840.625
$ws.Range("N121").Value = -6674086
# row 126
$ws.Range("H126").Value = 3701.5557
$ws.Range("I126").Value = 3092.5715
$ws.Range("J126").Value = 4089.0908
$ws.Range("K126").Value = 9277.7145
$ws.Range("L126").Value = 12267.2724
$ws.Range("M126").Value = -4337.7145
$ws.Range("N126").Value = -22147.2724
# row 133
$ws.Range("H133").Value = 6987.5
$ws.Range("I133").Value = 4316.6665
$ws.Range("J133").Value = 15000
$ws.Range("K133").Value = 12949.9995
$ws.Range("L133").Value = 45000
$ws.Range("M133").Value = -7889.999500000002
$ws.Range("N133").Value = -55120
# row 138
$ws.Range("H138").Value = 1800
$ws.Range("I138").Value = 1800
$ws.Range("K138").Value = 5400
$ws.Range("M138").Value = -260

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 4989.7676
$ws.Range("I102").Value = 5374.7427
$ws.Range("K102").Value = 5374.7427
$ws.Range("M102").Value = -3752.7427

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 405860.4
$ws.Range("I132").Value = 553253.5
$ws.Range("J132").Value = 7899
$ws.Range("K132").Value = 1659760.5
$ws.Range("L132").Value = 23697
$ws.Range("M132").Value = -1657230.5
$ws.Range("N132").Value = -28757

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 115
$ws.Range("H115").Value = 67400
$ws.Range("J115").Value = 67400
$ws.Range("L115").Value = 67400
$ws.Range("N115").Value = -70534
# row 132
$ws.Range("H132").Value = 21435.117
$ws.Range("I132").Value = 32841.95
$ws.Range("J132").Value = 6986.467
$ws.Range("K132").Value = 98525.84999999999
$ws.Range("L132").Value = 20959.401
$ws.Range("M132").Value = -95995.84999999999
$ws.Range("N132").Value = -26019.401
# row 136
$ws.Range("H136").Value = 258163.42
$ws.Range("I136").Value = 309184.8
$ws.Range("J136").Value = 3056.4
$ws.Range("K136").Value = 927554.3999999999
$ws.Range("L136").Value = 9169.200000000001
$ws.Range("M136").Value = -925004.3999999999
$ws.Range("N136").Value = -14269.2
